# Apply fix to "barriers" and "facilitators" rows in DescriptivesAll sheet.
# Commit message: "Fixed barriers and naming"
#
# Changes:
#   Row 13 (barriers):     F13 Mean 0.95 -> 0.9,  G13 SD 0.81 -> 0.75,
#                           H13 Range "0.00-   3.57" -> "0.00-   3.50"
#   Row 14 (facilitators): F14 Mean 1.03 -> 0.92, G14 SD 1.12 -> 1.02

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 13 - barriers
$ws.Range("F13").Value = 0.9
$ws.Range("G13").Value = 0.75
$ws.Range("H13").Value = "0.00-   3.50"

# Row 14 - facilitators
$ws.Range("F14").Value = 0.92
$ws.Range("G14").Value = 1.02
